$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Give the first paragraph a paragraph border (5-twip gap on every side)
# and widen its left indent from 120 to 225 twips (i.e. 6pt -> 11.25pt).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Rewrite the placeholder merge-field id carried by the first run.
$idRange = $p1.Range.Duplicate()
$found = $idRange.Find.Execute("**ID__AFFARS_5319_topic_10__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $idRange.Text = "**ID__AFFARS_5319_810_90__ID**"
}

# That id used to be followed by a run holding a single trailing space;
# remove it so the paragraph ends immediately after "**ID**".
$p1 = $d.Paragraphs(1)
$tailSpace = $p1.Range.Duplicate()
$tailSpace.Start = $tailSpace.End - 2
$tailSpace.End = $tailSpace.End - 1
if ($tailSpace.Text -eq " ") {
    $tailSpace.Delete()
}
